$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header text (G1/H1): rename "Primitive Diameter" to the more
# specific "Primitive Diameter/ Base Circle" and reorder so "Foot Diameter"
# comes first in the shared-string table.
$ws.Range("G1").Value = "Primitive Diameter/ Base Circle"
$ws.Range("H1").Value = "Foot Diameter"

# Update the existing sample row (row 2) with new gear measurements.
$ws.Range("A2").Value = 25
$ws.Range("B2").Value = 27

# Add a new gear measurement row (row 3) with the same formulas as row 2.
$ws.Range("A3").Value = 33
$ws.Range("B3").Value = 35
$ws.Range("C3").Formula = "=B3/(A3+2)"
$ws.Range("D3").Formula = "=1/C3"
$ws.Range("E3").Formula = "=A3/D3"
$ws.Range("F3").Formula = "=C3*PI()"
$ws.Range("G3").Formula = "=C3*A3"
$ws.Range("H3").Formula = "=G3-2*C3"

# Move the active selection to G3.
$ws.Range("G3").Select()
